$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A16").Borders.Item(5).LineStyle = 1
